$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2,4)
$c.Value = "'244.29"
$c.Style = "Normal"
$ws.Cells.Item(2,6).Value = '29-12-2022'
$c = $ws.Cells.Item(2,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(3,4)
$c.Value = "'23.95"
$c.Style = "Normal"
$ws.Cells.Item(3,6).Value = '29-12-2022'
$c = $ws.Cells.Item(3,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(4,4)
$c.Value = "'5.249"
$c.Style = "Normal"
$ws.Cells.Item(4,6).Value = '29-12-2022'
$c = $ws.Cells.Item(4,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(5,4)
$c.Value = "'0.05773"
$c.Style = "Normal"
$ws.Cells.Item(5,6).Value = '29-12-2022'
$c = $ws.Cells.Item(5,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(6,4)
$c.Value = "'6.480"
$c.Style = "Normal"
$ws.Cells.Item(6,6).Value = '29-12-2022'
$c = $ws.Cells.Item(6,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(7,4)
$c.Value = "'3.253"
$c.Style = "Normal"
$ws.Cells.Item(7,6).Value = '29-12-2022'
$c = $ws.Cells.Item(7,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(8,4)
$c.Value = "'0.8121"
$c.Style = "Normal"
$ws.Cells.Item(8,6).Value = '29-12-2022'
$c = $ws.Cells.Item(8,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(9,4)
$c.Value = "'0.8836"
$c.Style = "Normal"
$ws.Cells.Item(9,6).Value = '29-12-2022'
$c = $ws.Cells.Item(9,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(10,4)
$c.Value = "'0.1375"
$c.Style = "Normal"
$ws.Cells.Item(10,6).Value = '29-12-2022'
$c = $ws.Cells.Item(10,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(11,4)
$c.Value = "'0.06948"
$c.Style = "Normal"
$ws.Cells.Item(11,6).Value = '29-12-2022'
$c = $ws.Cells.Item(11,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(12,4)
$c.Value = "'0.03201"
$c.Style = "Normal"
$ws.Cells.Item(12,6).Value = '29-12-2022'
$c = $ws.Cells.Item(12,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(13,4)
$c.Value = "'0.03036"
$c.Style = "Normal"
$ws.Cells.Item(13,6).Value = '29-12-2022'
$c = $ws.Cells.Item(13,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(14,4)
$c.Value = "'0.09333"
$c.Style = "Normal"
$ws.Cells.Item(14,6).Value = '29-12-2022'
$c = $ws.Cells.Item(14,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(15,4)
$c.Value = "'3.819"
$c.Style = "Normal"
$ws.Cells.Item(15,6).Value = '29-12-2022'
$c = $ws.Cells.Item(15,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(16,4)
$c.Value = "'0.001520"
$c.Style = "Normal"
$ws.Cells.Item(16,6).Value = '29-12-2022'
$c = $ws.Cells.Item(16,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(17,4)
$c.Value = "'0.04701"
$c.Style = "Normal"
$ws.Cells.Item(17,6).Value = '29-12-2022'
$c = $ws.Cells.Item(17,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(18,4)
$c.Value = "'0.0006034"
$c.Style = "Normal"
$ws.Cells.Item(18,6).Value = '29-12-2022'
$c = $ws.Cells.Item(18,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(19,4)
$c.Value = "'0.006181"
$c.Style = "Normal"
$ws.Cells.Item(19,6).Value = '29-12-2022'
$c = $ws.Cells.Item(19,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(20,4)
$c.Value = "'0.001235"
$c.Style = "Normal"
$ws.Cells.Item(20,6).Value = '29-12-2022'
$c = $ws.Cells.Item(20,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(21,4)
$c.Value = "'0.004084"
$c.Style = "Normal"
$ws.Cells.Item(21,6).Value = '29-12-2022'
$c = $ws.Cells.Item(21,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(22,4)
$c.Value = "'0.00008691"
$c.Style = "Normal"
$ws.Cells.Item(22,6).Value = '29-12-2022'
$c = $ws.Cells.Item(22,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(23,4)
$c.Value = "'3.548"
$c.Style = "Normal"
$ws.Cells.Item(23,6).Value = '29-12-2022'
$c = $ws.Cells.Item(23,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(24,4)
$c.Value = "'2.144"
$c.Style = "Normal"
$ws.Cells.Item(24,6).Value = '29-12-2022'
$c = $ws.Cells.Item(24,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(25,4)
$c.Value = "'0.3165"
$c.Style = "Normal"
$ws.Cells.Item(25,6).Value = '29-12-2022'
$c = $ws.Cells.Item(25,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(26,6).Value = '29-12-2022'
$c = $ws.Cells.Item(26,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(27,6).Value = '29-12-2022'
$c = $ws.Cells.Item(27,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(28,6).Value = '29-12-2022'
$c = $ws.Cells.Item(28,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(29,6).Value = '29-12-2022'
$c = $ws.Cells.Item(29,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(30,6).Value = '29-12-2022'
$c = $ws.Cells.Item(30,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(31,6).Value = '29-12-2022'
$c = $ws.Cells.Item(31,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(32,6).Value = '29-12-2022'
$c = $ws.Cells.Item(32,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(33,6).Value = '29-12-2022'
$c = $ws.Cells.Item(33,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(34,6).Value = '29-12-2022'
$c = $ws.Cells.Item(34,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(35,6).Value = '29-12-2022'
$c = $ws.Cells.Item(35,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(36,6).Value = '29-12-2022'
$c = $ws.Cells.Item(36,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(37,6).Value = '29-12-2022'
$c = $ws.Cells.Item(37,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(38,6).Value = '29-12-2022'
$c = $ws.Cells.Item(38,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(39,6).Value = '29-12-2022'
$c = $ws.Cells.Item(39,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(40,4)
$c.Value = "'0.03724"
$c.Style = "Normal"
$ws.Cells.Item(40,6).Value = '29-12-2022'
$c = $ws.Cells.Item(40,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(41,2).Value = 'KickToken'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$c = $ws.Cells.Item(41,4)
$c.Value = "'0.006253"
$c.Style = "Normal"
$ws.Cells.Item(41,5).Value = '40KickTokenKICK'
$ws.Cells.Item(41,6).Value = '29-12-2022'
$c = $ws.Cells.Item(41,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(42,2).Value = 'BKEXToken'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$c = $ws.Cells.Item(42,4)
$c.Value = "'0.1050"
$c.Style = "Normal"
$ws.Cells.Item(42,5).Value = '41BKEXTokenBKK'
$ws.Cells.Item(42,6).Value = '29-12-2022'
$c = $ws.Cells.Item(42,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(43,2).Value = 'CEJI'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$c = $ws.Cells.Item(43,4)
$c.Value = "'0.002398"
$c.Style = "Normal"
$ws.Cells.Item(43,5).Value = '42CEJICEJIWorstin24h'
$ws.Cells.Item(43,6).Value = '29-12-2022'
$c = $ws.Cells.Item(43,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(44,4)
$c.Value = "'0.007894"
$c.Style = "Normal"
$ws.Cells.Item(44,6).Value = '29-12-2022'
$c = $ws.Cells.Item(44,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(45,4)
$c.Value = "'0.00005314"
$c.Style = "Normal"
$ws.Cells.Item(45,6).Value = '29-12-2022'
$c = $ws.Cells.Item(45,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(46,6).Value = '29-12-2022'
$c = $ws.Cells.Item(46,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(47,4)
$c.Value = "'0.4197"
$c.Style = "Normal"
$ws.Cells.Item(47,6).Value = '29-12-2022'
$c = $ws.Cells.Item(47,7)
$c.Value = "'0"
$c.Style = "Normal"

$c = $ws.Cells.Item(48,4)
$c.Value = "'0.002666"
$c.Style = "Normal"
$ws.Cells.Item(48,6).Value = '29-12-2022'
$c = $ws.Cells.Item(48,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(49,6).Value = '29-12-2022'
$c = $ws.Cells.Item(49,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(50,6).Value = '29-12-2022'
$c = $ws.Cells.Item(50,7)
$c.Value = "'0"
$c.Style = "Normal"

$ws.Cells.Item(51,6).Value = '29-12-2022'
$c = $ws.Cells.Item(51,7)
$c.Value = "'0"
$c.Style = "Normal"

Write-Output "Applied symbol list update for 29-12-2022"